# DE_table2_F2_dim10.xlsx update
# - Rename the "Gen" header/column (A) to "MaxFES" and replace its generation-count
#   values with the fractional MaxFES values.
# - Drop the "Run 50" run (its data column) so only Run 0..Run 49 remain.
# - Recompute the "Mean" column over the remaining 50 runs and place it directly
#   after "Run 49" (shifting left into what used to be the "Run 50" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: "Gen" -> "MaxFES" header + new values (rows 2-14) ---
$ws.Range("A1").Value = "MaxFES"

$maxFesValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxFesValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $maxFesValues[$i]
}

# --- Recompute Mean into column AZ (currently "Run 50"), then drop old BA ("Mean") ---
$ws.Range("AZ1").Value = "Mean"

$meanValues = @(
    11103477560.4985,
    3449748785.632158,
    530914226.1220242,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089,
    527597560.6226089
)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $meanValues[$i]
}

# Remove the old "Run 50" data / old "Mean" column entirely (column BA), so the
# new Mean values in AZ become the last column and the sheet shrinks to A:AZ.
$ws.Range("BA1:BA14").EntireColumn.Delete()
